# TemplateLeaveRequest.xlsx edit
# - "Main" sheet (leave request template): insert a new "Code" column (D)
#   before the "Tanggal" column. Downstream columns (Tanggal, Jumlah Hari,
#   Tipe Cuti, Keterangan Cuti / Alasan) shift one column to the right.
#   The "Tipe Cuti" helper comment + its data validation move from F to G.
# - "Example" sheet: same column insert, plus sample "Code" values for the
#   4 example rows (formatted in a distinct font), and refreshed sample
#   dates for the new permintaan-izin examples.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Main")
$wsExample = $wb.Worksheets.Item("Example")

# ---------------------------------------------------------------------
# 1) "Main" sheet — insert the Code column (D), header only.
# ---------------------------------------------------------------------

# Move the "Tipe Cuti" helper comment off column F before the insert so it
# isn't silently dropped, then re-create it one column over (G) afterwards.
$mainComment = $wsMain.Range("F1").Comment
$mainCommentText = $mainComment.Text()
$mainComment.Delete()

$wsMain.Columns.Item(4).Insert()

$wsMain.Range("D1").Value = "Code"
$wsMain.Columns.Item(4).ColumnWidth = 14.59

$newMainComment = $wsMain.Range("G1").AddComment($mainCommentText)
$newMainComment.Author = "ismail - [2010]"

# Data validation dropdown for "Tipe Cuti" now lives in column G.
$wsMain.Range("F2:F69").Validation.Delete()
$wsMain.Range("G2:G69").Validation.Add(3, 1, 1, '"Cuti, Cuti Bersama, Cuti Melahirkan, Cuti Menikah"')

$wsMain.Activate()
$wsMain.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) "Example" sheet — insert the Code column (D) + sample data.
# ---------------------------------------------------------------------

$exComment = $wsExample.Range("F1").Comment
$exCommentText = $exComment.Text()
$exComment.Delete()

$wsExample.Columns.Item(4).Insert()

$wsExample.Range("D1").Value = "Code"

$newExComment = $wsExample.Range("G1").AddComment($exCommentText)
$newExComment.Author = "ismail - [2010]"

$wsExample.Range("F2:F5").Validation.Delete()
$wsExample.Range("G2:G5").Validation.Add(3, 1, 1, '"Cuti, Cuti Bersama, Cuti Melahirkan, Cuti Menikah"')

# Sample request codes for the 4 example rows.
$wsExample.Range("D2").Value = "0774/SS/PERS-JKT/I/2025"
$wsExample.Range("D3").Value = "0774/SI/PERS-JKT/I/2025"
$wsExample.Range("D4").Value = "0774/SI/PERS-JKT/I/2025"
$wsExample.Range("D5").Value = "0774/SC/PERS-JKT/I/2025"

# Give the Code column its own look (bigger, dark-grey Segoe UI) and widen it.
$codeRange = $wsExample.Range("D2:D5")
$codeFont = $codeRange.Font
$codeFont.Name = "Segoe UI"
$codeFont.Size = 12
$codeFont.Color = 0x292521
$wsExample.Columns.Item(4).ColumnWidth = 27.76

# Refreshed example dates (moved from Tanggal column E, after the insert).
$wsExample.Range("E2").Value = 45667
$wsExample.Range("E3").Value = 45683
$wsExample.Range("E4").Value = 45299
$wsExample.Range("E5").Value = 45292

# Taller rows to match the bigger Code font used in the example rows.
$wsExample.Rows.Item(2).RowHeight = 17.25
$wsExample.Rows.Item(3).RowHeight = 17.25
$wsExample.Rows.Item(4).RowHeight = 17.25
$wsExample.Rows.Item(5).RowHeight = 17.25

$wsExample.Activate()
$wsExample.Range("F5").Select() | Out-Null

$wsMain.Activate()
